$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "primary track" settings values shared by every data row (columns
# B-J). They are stored as text (shared strings) in the workbook, not as
# numbers (F/SegClipping is the only genuine number).
$segHighStd = "21.3968"
$segLowStd = "3.6228"
$segBgFactor = "0.93868"
$segThreshold = "0.007254"
$segClipping = 1
$segWHMax = "0.035885"
$segWHMax2 = "7.08"
$segMinArea = "180.0114"
$segMinSumIntensity = "5.8158"

function Set-SettingsRow {
    param($rowIndex, $name)
    $ws.Range("A$rowIndex").Value = $name
    $ws.Range("B$rowIndex").Value = $segHighStd
    $ws.Range("C$rowIndex").Value = $segLowStd
    $ws.Range("D$rowIndex").Value = $segBgFactor
    $ws.Range("E$rowIndex").Value = $segThreshold
    $ws.Range("F$rowIndex").Value = $segClipping
    $ws.Range("G$rowIndex").Value = $segWHMax
    $ws.Range("H$rowIndex").Value = $segWHMax2
    $ws.Range("I$rowIndex").Value = $segMinArea
    $ws.Range("J$rowIndex").Value = $segMinSumIntensity
}

# Insert two new rows above the existing row 2 (Fluo-C2DL-MSC) for the new
# BF-C2DL-HSC and BF-C2DL-MuSC entries.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Insert a new row before the (now shifted to row 4) Fluo-C2DL-MSC row for
# DIC-C2DH-HeLa.
$ws.Rows.Item(4).Insert()

# Insert a new row before the (now shifted to row 14) PhC-C2DL-PSC row for
# PhC-C2DH-U373.
$ws.Rows.Item(13).Insert()

# Mark the text-valued columns (B-E, G-J) of every data row as Text before
# assigning, so that numeric-looking strings like "21.3968" are kept as
# text (matching the original workbook, where these are stored as shared
# strings, not numbers).
$textRanges = @($ws.Range("B2:E14"), $ws.Range("G2:J14"))
foreach ($r in $textRanges) { $r.NumberFormat = "@" }

# All 13 dataset rows now use the same updated "primary track" settings.
Set-SettingsRow 2 "BF-C2DL-HSC"
Set-SettingsRow 3 "BF-C2DL-MuSC"
Set-SettingsRow 4 "DIC-C2DH-HeLa"
Set-SettingsRow 5 "Fluo-C2DL-MSC"
Set-SettingsRow 6 "Fluo-C3DH-A549"
Set-SettingsRow 7 "Fluo-C3DH-H157"
Set-SettingsRow 8 "Fluo-C3DL-MDA231"
Set-SettingsRow 9 "Fluo-N2DH-GOWT1"
Set-SettingsRow 10 "Fluo-N2DL-HeLa"
Set-SettingsRow 11 "Fluo-N3DH-CE"
Set-SettingsRow 12 "Fluo-N3DH-CHO"
Set-SettingsRow 13 "PhC-C2DH-U373"
Set-SettingsRow 14 "PhC-C2DL-PSC"

# Restore the default style on those cells (keep them as text values, but
# drop the explicit text number format so the cell style matches the rest
# of the sheet).
foreach ($r in $textRanges) { $r.Style = "Normal" }
